$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# A7: "2018.07.05" looks like a valid calendar date, so a plain .Value
# assignment gets auto-converted to a date serial (and picks up a date
# number format). Route it through a text formula + paste-special-values
# so it lands as a genuine shared-string cell with no style change,
# matching how the existing Date column cells are stored.
$ws.Cells.Item($row, 1).Formula = "=""2018.07.05"""
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)

$ws.Cells.Item($row, 2).Value = "12:22:01"
$ws.Cells.Item($row, 3).Value = "RS"
$ws.Cells.Item($row, 4).Value = 32
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 100
$ws.Cells.Item($row, 7).Value = 100
$ws.Cells.Item($row, 8).Value = 250
$ws.Cells.Item($row, 9).Value = "N/A"
$ws.Cells.Item($row, 10).Value = 1
$ws.Cells.Item($row, 11).Value = "effective"
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 9.08
$ws.Cells.Item($row, 14).Value = 117
$ws.Cells.Item($row, 15).Value = 35.4
$ws.Cells.Item($row, 16).Value = 0.5413239918684565
